# Regenerate save_data column G ("K") with recomputed strike counts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 4
    3  = 7
    4  = 12
    5  = 6
    6  = 7
    7  = 8
    8  = 2
    9  = 5
    10 = 3
    11 = 9
    12 = 7
    13 = 11
    14 = 4
    15 = 6
    16 = 8
    17 = 9
    18 = 10
    19 = 8
    20 = 5
    21 = 2
    22 = 4
    23 = 5
    24 = 4
    25 = 10
    26 = 7
    27 = 3
    28 = 8
    29 = 6
    30 = 3
    31 = 7
    32 = 1
    33 = 6
    34 = 3
    35 = 3
    36 = 7
    37 = 3
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
